# Apply 2024-10-26 violent crime data updates across Citywide Totals,
# By Neighborhood, and individual neighborhood worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6645
$ws.Range("K3").Value = 6863
$ws.Range("G4").Value = 1491
$ws.Range("K4").Value = 1422
$ws.Range("K5").Value = 496
$ws.Range("K6").Value = 7551
$ws.Range("G7").Value = 24719
$ws.Range("K7").Value = 22977

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("K2").Value = 4
$ws.Range("K6").Value = 24

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 78
$ws.Range("K6").Value = 124
$ws.Range("K7").Value = 291

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 414
$ws.Range("K3").Value = 459
$ws.Range("K6").Value = 498
$ws.Range("K7").Value = 1503

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 172
$ws.Range("K7").Value = 496

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 252
$ws.Range("K6").Value = 314
$ws.Range("K7").Value = 996

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 128
$ws.Range("K6").Value = 88
$ws.Range("K7").Value = 378

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 223
$ws.Range("K3").Value = 257
$ws.Range("K7").Value = 779

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 131
$ws.Range("K6").Value = 195
$ws.Range("K7").Value = 537

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 100
$ws.Range("K3").Value = 161
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 389

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 199
$ws.Range("K7").Value = 694
$ws.Range("K8").Value = 1503
$ws.Range("K11").Value = 423
$ws.Range("K15").Value = 240
$ws.Range("K19").Value = 671
$ws.Range("K20").Value = 555
$ws.Range("K25").Value = 108
$ws.Range("K27").Value = 215
$ws.Range("K29").Value = 1244
$ws.Range("K33").Value = 996
$ws.Range("K34").Value = 129
$ws.Range("K36").Value = 292
$ws.Range("K37").Value = 779
$ws.Range("K38").Value = 24
$ws.Range("K41").Value = 159
$ws.Range("K42").Value = 847
$ws.Range("K44").Value = 191
$ws.Range("K47").Value = 154
$ws.Range("K48").Value = 292
$ws.Range("K51").Value = 288
$ws.Range("K52").Value = 609
$ws.Range("K53").Value = 291
$ws.Range("K54").Value = 454
$ws.Range("K57").Value = 86
$ws.Range("K58").Value = 16
$ws.Range("G63").Value = 293
$ws.Range("K63").Value = 60
$ws.Range("K64").Value = 143
$ws.Range("K65").Value = 537
$ws.Range("K67").Value = 897
$ws.Range("K73").Value = 208
$ws.Range("K78").Value = 261
$ws.Range("K79").Value = 574
$ws.Range("K80").Value = 83
$ws.Range("K82").Value = 26
$ws.Range("K83").Value = 496
$ws.Range("K85").Value = 1062
$ws.Range("K86").Value = 140
$ws.Range("K89").Value = 345
$ws.Range("K90").Value = 218
$ws.Range("K91").Value = 273
$ws.Range("K95").Value = 378
$ws.Range("K96").Value = 243
$ws.Range("K98").Value = 115
$ws.Range("K99").Value = 389
$ws.Range("G101").Value = 24719
$ws.Range("K101").Value = 22977

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 244
$ws.Range("K3").Value = 328
$ws.Range("K4").Value = 49
$ws.Range("K6").Value = 255
$ws.Range("K7").Value = 897

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K4").Value = 26
$ws.Range("K7").Value = 454

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 352
$ws.Range("K6").Value = 364
$ws.Range("K7").Value = 1244

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 68
$ws.Range("K6").Value = 139
$ws.Range("K7").Value = 292

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 201
$ws.Range("K6").Value = 223
$ws.Range("K7").Value = 671

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K3").Value = 50
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 159

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 258
$ws.Range("K6").Value = 314
$ws.Range("K7").Value = 847

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 77
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 243

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 273

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 185
$ws.Range("K6").Value = 142
$ws.Range("K7").Value = 574

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 193
$ws.Range("K4").Value = 26
$ws.Range("K7").Value = 555

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 88
$ws.Range("K7").Value = 292

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 227
$ws.Range("K4").Value = 24
$ws.Range("K6").Value = 191
$ws.Range("K7").Value = 694

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 45
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 88
$ws.Range("K3").Value = 61
$ws.Range("K4").Value = 15
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 110
$ws.Range("K7").Value = 423

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 107
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 345

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 81
$ws.Range("K7").Value = 218

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 96
$ws.Range("K7").Value = 288

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 349
$ws.Range("K3").Value = 370
$ws.Range("K4").Value = 54
$ws.Range("K7").Value = 1062

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K5").Value = 14
$ws.Range("K6").Value = 26

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K2").Value = 19
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 172
$ws.Range("K7").Value = 609

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("K2").Value = 2
$ws.Range("K7").Value = 16

Write-Output "Updated 173 cells across 46 worksheets for 2024-10-26."